$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new headers I1 ("I0") and J1 ("IF"), matching the formatting of
# --- the existing header row (bold, bordered, centered) by copying the
# --- format from the neighboring header cell H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Fill in the new I (I0) and J (IF) data columns for rows 2-64.
$data = @(
    @(7,7),
    @(6,6),
    @(7,7),
    @(7,7),
    @(5,5),
    @(5,5),
    @(8,8),
    @(6,6),
    @(7,7),
    @(5,5),
    @(6,6),
    @(7,7),
    @(6,6),
    @(7,7),
    @(6,6),
    @(7,7),
    @(8,8),
    @(5,5),
    @(6,7),
    @(7,7),
    @(7,7),
    @(6,6),
    @(6,6),
    @(7,7),
    @(5,6),
    @(8,8),
    @(6,6),
    @(6,6),
    @(5,5),
    @(6,6),
    @(5,5),
    @(7,7),
    @(6,6),
    @(6,7),
    @(8,8),
    @(6,6),
    @(7,7),
    @(4,5),
    @(4,4),
    @(7,7),
    @(6,6),
    @(5,5),
    @(6,6),
    @(7,7),
    @(6,6),
    @(8,8),
    @(9,9),
    @(6,6),
    @(7,7),
    @(4,5),
    @(6,6),
    @(5,6),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(5,6),
    @(5,6),
    @(7,7),
    @(8,8),
    @(6,6),
    @(7,7),
    @(4,4)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = 2 + $idx
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
